$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ND-AD-dAD")

# Update existing values in the first table (Dining Philosophers)
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2

# Fill in the newly-populated rows in the second table (German Cache Coherence Protocol)
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2

$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1

$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1

$ws.Range("B19").Value = 2
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 1

# Update the selected cell shown in the sheet view
$ws.Activate()
$ws.Range("C23").Select()
